$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M25").Value = 6325.86
$ws1.Range("M26").Value = "2 de 24"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F25").Value = 6325.86
$ws2.Range("F26").Value = 5129.469999999999

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D12").Value = 6521.5
$ws3.Range("E12").Value = 21433.48
$ws3.Range("F12").Value = 0.2332858045328596
$ws3.Range("D14").Value = 5129.47
$ws3.Range("E14").Value = 37073.91110009469
$ws3.Range("F14").Value = 0.1215416837772861
$ws3.Columns.Item(6).ColumnWidth = 24.166666666666668
